# Auto-generated edit script: updates Titan_Profits workbook per commit diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 1235306.1
$ws.Range("I107").Value = 1587749.9
$ws.Range("J107").Value = 1753
$ws.Range("K107").Value = 1587749.9
$ws.Range("L107").Value = 1753
$ws.Range("M107").Value = -1585829.9
$ws.Range("N107").Value = -5593
$ws.Range("H116").Value = 4196126
$ws.Range("I116").Value = 11531899
$ws.Range("J116").Value = 4255.905
$ws.Range("K116").Value = 11531899
$ws.Range("L116").Value = 4255.905
$ws.Range("M116").Value = -11528457
$ws.Range("N116").Value = -11139.905
$ws.Range("H125").Value = 18685818
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 18685818
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 168172362
$ws.Range("N125").Value = -168177282
$ws.Range("H132").Value = 226794.62
$ws.Range("I132").Value = 283314.12
$ws.Range("J132").Value = 39845.54
$ws.Range("K132").Value = 849942.36
$ws.Range("L132").Value = 119536.62
$ws.Range("M132").Value = -847412.36
$ws.Range("H136").Value = 43500
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 43500
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 43500
$ws.Range("N136").Value = -53700

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 26151.979
$ws.Range("I32").Value = 2865.0293
$ws.Range("J32").Value = 98129.82000000001
$ws.Range("K32").Value = 2865.0293
$ws.Range("L32").Value = 98129.82000000001
$ws.Range("M32").Value = -2578.0293
$ws.Range("N32").Value = -98703.82000000001
$ws.Range("H74").Value = 7901.421
$ws.Range("I74").Value = 1688.3077
$ws.Range("J74").Value = 21363.166
$ws.Range("K74").Value = 1688.3077
$ws.Range("L74").Value = 21363.166
$ws.Range("M74").Value = -814.3077000000001
$ws.Range("N74").Value = -23111.166
$ws.Range("H77").Value = 7901.421
$ws.Range("I77").Value = 1688.3077
$ws.Range("J77").Value = 21363.166
$ws.Range("K77").Value = 8441.538500000001
$ws.Range("L77").Value = 106815.83
$ws.Range("M77").Value = -4073.538500000001
$ws.Range("N77").Value = -115551.83
$ws.Range("H122").Value = 1900.3334
$ws.Range("I122").Value = 1800
$ws.Range("J122").Value = 2101
$ws.Range("K122").Value = 5400
$ws.Range("L122").Value = 6303
$ws.Range("M122").Value = -2950
$ws.Range("H139").Value = 54500
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 54500
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 54500
$ws.Range("N139").Value = -64780

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1465.5758
$ws.Range("I20").Value = 1294.5652
$ws.Range("J20").Value = 1858.9
$ws.Range("K20").Value = 1294.5652
$ws.Range("L20").Value = 1858.9
$ws.Range("M20").Value = -1047.5652
$ws.Range("N20").Value = -2352.9
$ws.Range("H94").Value = 793
$ws.Range("I94").Value = 793
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 793
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -342
$ws.Range("N94").ClearContents()
$ws.Range("H105").Value = 3100.0789
$ws.Range("I105").Value = 3310.6843
$ws.Range("J105").Value = 2889.4736
$ws.Range("K105").Value = 3310.6843
$ws.Range("L105").Value = 2889.4736
$ws.Range("M105").Value = -1563.6843
$ws.Range("N105").Value = -6383.473599999999
$ws.Range("H134").Value = 7204.6665
$ws.Range("I134").Value = 6500
$ws.Range("J134").Value = 7406
$ws.Range("K134").Value = 19500
$ws.Range("L134").Value = 22218
$ws.Range("M134").Value = -16965
$ws.Range("N134").Value = -27288

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 200120200
$ws.Range("I4").Value = 500000
$ws.Range("J4").Value = 250025250
$ws.Range("K4").Value = 500000
$ws.Range("L4").Value = 250025250
$ws.Range("M4").Value = -499888
$ws.Range("N4").Value = -250025474
$ws.Range("H31").Value = 1745.8077
$ws.Range("I31").Value = 1320.8125
$ws.Range("J31").Value = 2425.8
$ws.Range("K31").Value = 1320.8125
$ws.Range("L31").Value = 2425.8
$ws.Range("M31").Value = -1025.8125
$ws.Range("N31").Value = -3015.8
$ws.Range("H34").Value = 1745.8077
$ws.Range("I34").Value = 1320.8125
$ws.Range("J34").Value = 2425.8
$ws.Range("K34").Value = 1320.8125
$ws.Range("L34").Value = 2425.8
$ws.Range("M34").Value = -1118.8125
$ws.Range("N34").Value = -2829.8
$ws.Range("H105").Value = 488.75
$ws.Range("I105").Value = 488.75
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 488.75
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 1258.25
$ws.Range("N105").ClearContents()
$ws.Range("H107").Value = 406.9091
$ws.Range("I107").Value = 154
$ws.Range("J107").Value = 710.4
$ws.Range("K107").Value = 154
$ws.Range("L107").Value = 710.4
$ws.Range("M107").Value = 1766
$ws.Range("N107").Value = -4550.4
$ws.Range("H134").Value = 2760.75
$ws.Range("I134").Value = 1303.6842
$ws.Range("J134").Value = 5836.778
$ws.Range("K134").Value = 3911.0526
$ws.Range("L134").Value = 17510.334
$ws.Range("M134").Value = -1376.0526
$ws.Range("N134").Value = -22580.334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 786
$ws.Range("I122").Value = 540.6
$ws.Range("J122").Value = 1399.5
$ws.Range("K122").Value = 4865.400000000001
$ws.Range("L122").Value = 12595.5
$ws.Range("M122").Value = -2415.400000000001
$ws.Range("H139").Value = 2030.8
$ws.Range("I139").Value = 1301.6666
$ws.Range("J139").Value = 4947.3335
$ws.Range("K139").Value = 3904.9998
$ws.Range("L139").Value = 14842.0005
$ws.Range("M139").Value = 1235.0002
$ws.Range("N139").Value = -25122.0005
$ws.Range("H140").Value = 8407.267
$ws.Range("I140").Value = 11764.333
$ws.Range("J140").Value = 3371.6667
$ws.Range("K140").Value = 35292.999
$ws.Range("L140").Value = 10115.0001
$ws.Range("M140").Value = -30112.999
$ws.Range("N140").Value = -20475.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2418.1516
$ws.Range("I80").Value = 2388.8518
$ws.Range("J80").Value = 2550
$ws.Range("K80").Value = 2388.8518
$ws.Range("L80").Value = 2550
$ws.Range("M80").Value = -1390.8518
$ws.Range("N80").Value = -4546
$ws.Range("H83").Value = 2418.1516
$ws.Range("I83").Value = 2388.8518
$ws.Range("J83").Value = 2550
$ws.Range("K83").Value = 11944.259
$ws.Range("L83").Value = 12750
$ws.Range("M83").Value = -6952.259
$ws.Range("N83").Value = -22734
$ws.Range("H122").Value = 2224323.5
$ws.Range("I122").Value = 3705372.8
$ws.Range("J122").Value = 2750
$ws.Range("K122").Value = 11116118.4
$ws.Range("L122").Value = 8250
$ws.Range("M122").Value = -11113668.4
$ws.Range("N122").Value = -13150
$ws.Range("H132").Value = 3058.1738
$ws.Range("I132").Value = 1899
$ws.Range("J132").Value = 4861.3335
$ws.Range("K132").Value = 5697
$ws.Range("L132").Value = 14584.0005
$ws.Range("M132").Value = -3167
$ws.Range("H138").Value = 54900
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 54900
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 54900
$ws.Range("N138").Value = -65180

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3404.762
$ws.Range("I7").Value = 2666.6667
$ws.Range("J7").Value = 3527.7778
$ws.Range("K7").Value = 2666.6667
$ws.Range("L7").Value = 3527.7778
$ws.Range("M7").Value = -2554.6667
$ws.Range("N7").Value = -3751.7778
$ws.Range("H82").Value = 33087.5
$ws.Range("I82").Value = 43391.668
$ws.Range("J82").Value = 2175
$ws.Range("K82").Value = 43391.668
$ws.Range("L82").Value = 2175
$ws.Range("M82").Value = -43030.668
$ws.Range("N82").Value = -2897
$ws.Range("H85").Value = 33087.5
$ws.Range("I85").Value = 43391.668
$ws.Range("J85").Value = 2175
$ws.Range("K85").Value = 43391.668
$ws.Range("L85").Value = 2175
$ws.Range("M85").Value = -42143.668
$ws.Range("N85").Value = -4671
$ws.Range("H122").Value = 3800
$ws.Range("I122").Value = 2500
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 7500
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -5050
$ws.Range("N122").Value = -16900
$ws.Range("H126").Value = 3404.762
$ws.Range("I126").Value = 2666.6667
$ws.Range("J126").Value = 3527.7778
$ws.Range("K126").Value = 8000.000100000001
$ws.Range("L126").Value = 10583.3334
$ws.Range("M126").Value = -5530.000100000001
$ws.Range("N126").Value = -15523.3334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 101169.7
$ws.Range("I122").Value = 112299.664
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 336898.992
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -334448.992
$ws.Range("N122").Value = -7900
$ws.Range("H123").Value = 35617.332
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 35617.332
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 35617.332
$ws.Range("N123").Value = -45417.332
$ws.Range("H126").Value = 38173.184
$ws.Range("I126").Value = 53625.21
$ws.Range("J126").Value = 1474.625
$ws.Range("K126").Value = 160875.63
$ws.Range("L126").Value = 4423.875
$ws.Range("M126").Value = -158405.63
$ws.Range("N126").Value = -9363.875
$ws.Range("H136").Value = 23883556
$ws.Range("I136").Value = 41792880
$ws.Range("J136").Value = 4457
$ws.Range("K136").Value = 125378640
$ws.Range("L136").Value = 13371
$ws.Range("M136").Value = -125376090
$ws.Range("N136").Value = -18471

Write-Host "Applied edits: sets and clears complete"